$wb = $excel.ActiveWorkbook

# --- Update the TRP sheet's saved selection (A1:I53 -> A1:E53) ---
$trp = $wb.Worksheets.Item("TRP")
$null = $trp.Range("A1:E53").Select()

# --- Add the new "TRP+H2O" sheet after the last sheet ("TRP+Aux") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TRP+H2O"

# Copy the profile data (columns A-E, rows 1-53) from TRP into the new sheet,
# preserving number formats/styles. Row 1 only has two populated cells
# (A1, B1) so it is copied separately to avoid spilling blank cells into
# C1:E1.
$null = $trp.Range("A1:B1").Copy($newSheet.Range("A1"))
$null = $trp.Range("A2:E53").Copy($newSheet.Range("A2"))

# The new sheet's header (B1) reads "TRP-H2O" rather than "TRP" (new shared string).
$newSheet.Range("B1").Value = "TRP-H2O"

# Make the new sheet the active tab with B2 selected.
$null = $newSheet.Activate()
$null = $newSheet.Range("B2").Select()
